$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 100.4
$ws.Range("I4").Value = 125
$ws.Range("J4").Value = 51.2
$ws.Range("K4").Value = 125
$ws.Range("L4").Value = 51.2
$ws.Range("M4").Value = -11
$ws.Range("N4").Value = -279.2
$ws.Range("H42").Value = 2642.1667
$ws.Range("I42").Value = 212.14285
$ws.Range("J42").Value = 6044.2
$ws.Range("K42").Value = 636.4285500000001
$ws.Range("L42").Value = 18132.6
$ws.Range("M42").Value = -406.4285500000001
$ws.Range("N42").Value = -18592.6
$ws.Range("H62").Value = 3843.4285
$ws.Range("J62").Value = 3212.8333
$ws.Range("L62").Value = 3212.8333
$ws.Range("N62").Value = -4460.8333
$ws.Range("H65").Value = 3843.4285
$ws.Range("J65").Value = 3212.8333
$ws.Range("L65").Value = 16064.1665
$ws.Range("N65").Value = -22304.1665
$ws.Range("H100").Value = 5002889
$ws.Range("I100").Value = 2522.4546
$ws.Range("J100").Value = 8931748
$ws.Range("K100").Value = 2522.4546
$ws.Range("L100").Value = 8931748
$ws.Range("M100").Value = -1981.4546
$ws.Range("N100").Value = -8932830
$ws.Range("H116").Value = 3508.3062
$ws.Range("I116").Value = 3327.2903
$ws.Range("J116").Value = 3820.0557
$ws.Range("K116").Value = 3327.2903
$ws.Range("L116").Value = 3820.0557
$ws.Range("M116").Value = 114.7096999999999
$ws.Range("N116").Value = -10704.0557

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1073.5358
$ws.Range("I45").Value = 955.1818
$ws.Range("K45").Value = 955.1818
$ws.Range("M45").Value = -578.1818
$ws.Range("H61").Value = 2893.6584
$ws.Range("I61").Value = 2726.4
$ws.Range("K61").Value = 2726.4
$ws.Range("M61").Value = -2514.4
$ws.Range("H74").Value = 2289.2646
$ws.Range("I74").Value = 1963.88
$ws.Range("K74").Value = 1963.88
$ws.Range("M74").Value = -1089.88
$ws.Range("H77").Value = 2289.2646
$ws.Range("I77").Value = 1963.88
$ws.Range("K77").Value = 9819.400000000001
$ws.Range("M77").Value = -5451.400000000001
$ws.Range("H102").Value = 33368476
$ws.Range("I102").Value = 45456204
$ws.Range("K102").Value = 45456204
$ws.Range("M102").Value = -45454582
$ws.Range("H110").Value = 6897817
$ws.Range("I110").Value = 8001195
$ws.Range("K110").Value = 8001195
$ws.Range("M110").Value = -7999150
$ws.Range("H132").Value = 6074.34
$ws.Range("I132").Value = 4073.447
$ws.Range("J132").Value = 17412.732
$ws.Range("K132").Value = 12220.341
$ws.Range("L132").Value = 52238.196
$ws.Range("M132").Value = -9690.341
$ws.Range("N132").Value = -57298.196
$ws.Range("H136").Value = 2893.6584
$ws.Range("I136").Value = 2726.4
$ws.Range("K136").Value = 8179.200000000001
$ws.Range("M136").Value = -5629.200000000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2624.647
$ws.Range("I86").Value = 2741.3635
$ws.Range("J86").Value = 2410.6667
$ws.Range("K86").Value = 2741.3635
$ws.Range("L86").Value = 2410.6667
$ws.Range("M86").Value = -1618.3635
$ws.Range("N86").Value = -4656.6667
$ws.Range("H89").Value = 2624.647
$ws.Range("I89").Value = 2741.3635
$ws.Range("J89").Value = 2410.6667
$ws.Range("K89").Value = 13706.8175
$ws.Range("L89").Value = 12053.3335
$ws.Range("M89").Value = -8090.817499999999
$ws.Range("N89").Value = -23285.3335
$ws.Range("H94").Value = 27779788
$ws.Range("I94").Value = 1101.8276
$ws.Range("K94").Value = 1101.8276
$ws.Range("M94").Value = -650.8276000000001
$ws.Range("H99").Value = 4242
$ws.Range("I99").Value = 4115.8335
$ws.Range("K99").Value = 4115.8335
$ws.Range("M99").Value = -2617.8335
$ws.Range("H105").Value = 2514.8333
$ws.Range("I105").Value = 2072.5
$ws.Range("K105").Value = 2072.5
$ws.Range("M105").Value = -325.5
$ws.Range("H134").Value = 50003576
$ws.Range("I134").Value = 62504084
$ws.Range("K134").Value = 187512252
$ws.Range("M134").Value = -187509717

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1813
$ws.Range("J31").Value = 1820.2
$ws.Range("L31").Value = 1820.2
$ws.Range("N31").Value = -2410.2
$ws.Range("H34").Value = 1813
$ws.Range("J34").Value = 1820.2
$ws.Range("L34").Value = 1820.2
$ws.Range("N34").Value = -2224.2
$ws.Range("H134").Value = 1370.174
$ws.Range("I134").Value = 1281.75
$ws.Range("K134").Value = 3845.25
$ws.Range("M134").Value = -1310.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 125.411766
$ws.Range("I12").Value = 128.33333
$ws.Range("J12").Value = 123.818184
$ws.Range("K12").Value = 384.99999
$ws.Range("L12").Value = 371.454552
$ws.Range("M12").Value = -211.99999
$ws.Range("N12").Value = -717.454552
$ws.Range("H99").Value = 4015
$ws.Range("I99").Value = 4015
$ws.Range("K99").Value = 12045
$ws.Range("M99").Value = -9799
$ws.Range("H134").Value = 3785.5715
$ws.Range("I134").Value = 2649.8333
$ws.Range("K134").Value = 7949.499899999999
$ws.Range("M134").Value = -2879.499899999999
$ws.Range("H139").Value = 11860.538
$ws.Range("I139").Value = 4563
$ws.Range("K139").Value = 13689
$ws.Range("M139").Value = -8549
$ws.Range("H140").Value = 1123.2858
$ws.Range("I140").Value = 1123.2858
$ws.Range("K140").Value = 3369.8574
$ws.Range("M140").Value = 1810.1426

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28134.465
$ws.Range("I70").Value = 64776.223
$ws.Range("J70").Value = 10777.842
$ws.Range("K70").Value = 64776.223
$ws.Range("L70").Value = 10777.842
$ws.Range("M70").Value = -64506.223
$ws.Range("N70").Value = -11317.842
$ws.Range("H73").Value = 28134.465
$ws.Range("I73").Value = 64776.223
$ws.Range("J73").Value = 10777.842
$ws.Range("K73").Value = 64776.223
$ws.Range("L73").Value = 10777.842
$ws.Range("M73").Value = -63840.223
$ws.Range("N73").Value = -12649.842
$ws.Range("H80").Value = 6840.25
$ws.Range("I80").Value = 4050.8333
$ws.Range("K80").Value = 4050.8333
$ws.Range("M80").Value = -3052.8333
$ws.Range("H83").Value = 6840.25
$ws.Range("I83").Value = 4050.8333
$ws.Range("K83").Value = 20254.1665
$ws.Range("M83").Value = -15262.1665
$ws.Range("H97").Value = 621.7222
$ws.Range("I97").Value = 581.13336
$ws.Range("K97").Value = 581.13336
$ws.Range("M97").Value = -85.13336000000004
$ws.Range("H113").Value = 2047.7368
$ws.Range("I113").Value = 1842.5
$ws.Range("J113").Value = 2399.5715
$ws.Range("K113").Value = 1842.5
$ws.Range("L113").Value = 2399.5715
$ws.Range("M113").Value = 327.5
$ws.Range("N113").Value = -6739.5715
$ws.Range("H126").Value = 5968.385
$ws.Range("I126").Value = 5242.4287
$ws.Range("J126").Value = 6815.3335
$ws.Range("K126").Value = 15727.2861
$ws.Range("L126").Value = 20446.0005
$ws.Range("M126").Value = -13257.2861
$ws.Range("N126").Value = -25386.0005
$ws.Range("H132").Value = 2035.5897
$ws.Range("I132").Value = 1555.8928
$ws.Range("J132").Value = 3256.6365
$ws.Range("K132").Value = 4667.678400000001
$ws.Range("L132").Value = 9769.9095
$ws.Range("M132").Value = -2137.678400000001
$ws.Range("N132").Value = -14829.9095

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1265.9231
$ws.Range("I16").Value = 1134.3549
$ws.Range("K16").Value = 1134.3549
$ws.Range("M16").Value = -964.3549
$ws.Range("H93").Value = 3772.5715
$ws.Range("I93").Value = 2393.2856
$ws.Range("J93").Value = 6531.143
$ws.Range("K93").Value = 2393.2856
$ws.Range("L93").Value = 6531.143
$ws.Range("M93").Value = -1145.2856
$ws.Range("N93").Value = -9027.143
$ws.Range("H100").Value = 3937
$ws.Range("J100").Value = 5732.6665
$ws.Range("L100").Value = 5732.6665
$ws.Range("N100").Value = -6814.6665
$ws.Range("H132").Value = 4252.623
$ws.Range("I132").Value = 4146.2764
$ws.Range("J132").Value = 4609.643
$ws.Range("K132").Value = 12438.8292
$ws.Range("L132").Value = 13828.929
$ws.Range("M132").Value = -9908.8292
$ws.Range("N132").Value = -18888.929

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13160748
$ws.Range("I122").Value = 17243744
$ws.Range("K122").Value = 51731232
$ws.Range("M122").Value = -51728782
$ws.Range("H126").Value = 17545294
$ws.Range("I126").Value = 23810886
$ws.Range("K126").Value = 71432658
$ws.Range("M126").Value = -71430188
